# Demand module: add LES-CES function for household.
# - Rename "elasFU" -> "elasFU_CES"
# - Insert a new sheet "elasFU_LES" right after "elasFU_CES" (before "elasTRADE")
# - Populate elasFU_LES with the LES calibration/demand data
# - Make elasFU_LES the active tab, restore a saved selection on elasTRADE

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename elasFU -> elasFU_CES
# ---------------------------------------------------------------------------
$cesSheet = $wb.Worksheets.Item("elasFU")
$cesSheet.Name = "elasFU_CES"

# ---------------------------------------------------------------------------
# 2. Insert the new elasFU_LES sheet right after elasFU_CES
# ---------------------------------------------------------------------------
$lesSheet = $wb.Worksheets.Add($null, $cesSheet)
$lesSheet.Name = "elasFU_LES"

# Sector codes (column A, rows 2-36) - reused labels already present in the
# workbook (same order as the elasTRADE / elasKL sheets).
$sectors = @(
    "pPARI", "pWHEA", "pOCER", "pFVEG", "pOILS", "pSUGB", "pFIBR", "pOTHC",
    "pANIM", "pFORE", "pFISH", "pFOSM", "pOTHM", "pFBTO", "pTXWO", "pCOKE",
    "pREFN", "pCHEM", "pRUBP", "pNMMP", "pMETP", "pELEC", "pMACH", "pELCF",
    "pELCG", "pTRDI", "pHWAT", "pWATR", "pCONS", "pTRAD", "pHORE", "pTRAN",
    "pREBA", "pPUBO", "pWAST"
)

# elasINC (column B) values per row (rows 2-36, i.e. pPARI..pWAST)
$elasINC = @(
    0.32, 0.32, 0.32, 0.32, 0.32, 0.32, 0.32, 0.32, 0.32, 0.32, 0.32, 0.32,
    0.32, 0.32, 0.82, 0.97, 0.97, 1.29, 1.29, 1.29, 1.29, 1.04, 1.04, 0.97,
    0.97, 0.97, 0.97, 0.97, 1.29, 1.29, 1.29, 1.23, 1.29, 1.18, 1.29
)

# frisch (column C) is constant -1.54 for every row
$frisch = -1.54

# Row 1 headers
$lesSheet.Range("B1").Value = "elasINC"
$lesSheet.Range("C1").Value = "frisch"

# Data rows 2-36
for ($i = 0; $i -lt $sectors.Length; $i++) {
    $row = $i + 2
    $lesSheet.Cells.Item($row, 1).Value = $sectors[$i]
    $lesSheet.Cells.Item($row, 2).Value = $elasINC[$i]
    $lesSheet.Cells.Item($row, 3).Value = $frisch
}

# Apply the quote-prefixed label style (same style as the sector-name column
# used throughout the rest of the workbook) to A1:A36 in one shot.
$styleSrc = $wb.Worksheets.Item("elasTRADE").Range("A2")
$styleSrc.Copy()
$lesSheet.Range("A1:A36").PasteSpecial(-4122)
$lesSheet.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Source/citation notes in column G
#    (written in this exact order so new shared-string indices line up with
#    the saved workbook: elasINC, frisch, G3, G5, G7, G1, G4)
# ---------------------------------------------------------------------------
$g3 = $lesSheet.Range("G3")
$g3.Value = "Arjan Lejour & Paul Veenendaal & Gerard Verweij & Nico van Leeuwen, 2006."
$g3.Font.Name = "Arial Unicode MS"
$g3.Font.Size = 10
$g3.VerticalAlignment = -4108

$g5 = $lesSheet.Range("G5")
$g5.Value = "CPB Document 111, CPB Netherlands Bureau for Economic Policy Analysis."
$g5.Font.Name = "Arial Unicode MS"
$g5.Font.Size = 10
$g5.VerticalAlignment = -4108

$g7 = $lesSheet.Range("G7")
$g7.Value = "Original source: GTAP5/6, Dimaranan and McDougall (2002, 2005)."
$g7.Font.Name = "Arial Unicode MS"
$g7.Font.Size = 10
$g7.Font.Color = 0
$g7.VerticalAlignment = -4108

$lesSheet.Range("G1").Value = "The values are taken and assigned to the model sector from WorldScan description, for OECD countries (page 66)"

$g4 = $lesSheet.Range("G4")
$g4.Value = """Worldscan: a model for international economic policy analysis"""
$g4.Font.Name = "Arial Unicode MS"
$g4.Font.Size = 10
$g4.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Restore the remembered selection on elasTRADE, then make elasFU_LES the
#    active / selected tab (matches activeTab="1" in the saved workbook).
# ---------------------------------------------------------------------------
$tradeSheet = $wb.Worksheets.Item("elasTRADE")
$tradeSheet.Range("A2:A36").Select()

$lesSheet.Activate()
$lesSheet.Range("A1").Select()

Write-Output "done"
